$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '69.712.90'
$ws.Range('E2').Value = '  +0.40%  '
$ws.Range('D3').Value = '3.704.03'
$ws.Range('E3').Value = '  +0.38%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').Value = "'677.78"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.15%  '
$ws.Range('D6').Value = "'161.64"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.70%  '
$ws.Range('E7').Value = '  +0.09%  '
$ws.Range('D8').Value = "'0.496"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.39%  '
$ws.Range('E9').Value = '  +1.42%  '
$ws.Range('D10').Value = "'7.15"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.61%  '
$ws.Range('E11').Value = '  +1.86%  '
$ws.Range('E12').Value = '  +0.60%  '
$ws.Range('D13').Value = "'32.80"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.76%  '
$ws.Range('D14').Value = '3.696.84'
$ws.Range('E14').Value = '  +0.28%  '
$ws.Range('D15').Value = '69.741.90'
$ws.Range('E15').Value = '  +0.46%  '
$ws.Range('E16').Value = '  +1.80%  '
$ws.Range('E17').Value = '  +1.29%  '
$ws.Range('D18').Value = "'6.50"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +1.17%  '
$ws.Range('D19').Value = "'473.45"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.55%  '
$ws.Range('D20').Value = "'9.80"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -2.19%  '
$ws.Range('E21').Value = '  +0.45%  '
$ws.Range('D22').Value = "'80.49"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.76%  '
$ws.Range('D23').Value = '3.852.25'
$ws.Range('E23').Value = '  +0.38%  '
$ws.Range('E24').Value = '  +0.08%  '
$ws.Range('E26').Value = '  -0.11%  '
$ws.Range('D27').Value = "'9.12"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E29').Value = '  +0.62%  '
$ws.Range('E30').Value = '  +1.07%  '
$ws.Range('E31').Value = '  -0.06%  '
$ws.Range('E32').Value = '  +0.13%  '
$ws.Range('B33').Value = 'Kaspa'
$ws.Range('C33').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D33').Value = "'0.167"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +4.30%  '
$ws.Range('B34').Value = 'EthereumClassic'
$ws.Range('C34').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D34').Value = "'26.96"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.31%  '
$ws.Range('D35').Value = '3.695.37'
$ws.Range('E35').Value = '  +0.87%  '
$ws.Range('E36').Value = '  +3.86%  '
$ws.Range('D37').Value = "'6.22"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.83%  '
$ws.Range('E39').Value = '  -0.12%  '
$ws.Range('E40').Value = '  +0.02%  '
$ws.Range('D41').Value = "'0.0905"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.32%  '
$ws.Range('D42').Value = "'0.946"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.23%  '
$ws.Range('D43').Value = "'166.81"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.41%  '
$ws.Range('E44').Value = '  -0.99%  '
$ws.Range('E45').Value = '  +2.04%  '
$ws.Range('D46').Value = "'28.32"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.17%  '
$ws.Range('D47').Value = "'0.000280"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.28%  '
$ws.Range('E48').Value = '  -1.76%  '
$ws.Range('E49').Value = '  -0.34%  '
$ws.Range('D50').Value = "'7.90"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +1.22%  '
$ws.Range('E51').Value = '  +1.67%  '
